{"js": "// Replace each two-digit-by-two-digit multiplication expression in the\n// document with its new value. Each original expression is unique within\n// the document, so searching for the exact old text and replacing it is\n// unambiguous.\nconst pairs = [\n  [\"26\u00d725=\", \"72\u00d721=\"],\n  [\"92\u00d798=\", \"62\u00d722=\"],\n  [\"79\u00d736=\", \"74\u00d748=\"],\n  [\"16\u00d776=\", \"33\u00d747=\"],\n  [\"50\u00d793=\", \"65\u00d753=\"],\n  [\"39\u00d766=\", \"79\u00d754=\"],\n  [\"69\u00d788=\", \"21\u00d790=\"],\n  [\"14\u00d725=\", \"22\u00d787=\"],\n  [\"56\u00d779=\", \"64\u00d774=\"],\n  [\"53\u00d713=\", \"58\u00d798=\"],\n  [\"14\u00d786=\", \"61\u00d739=\"],\n  [\"92\u00d775=\", \"94\u00d758=\"],\n  [\"66\u00d723=\", \"87\u00d728=\"],\n  [\"80\u00d742=\", \"36\u00d746=\"],\n  [\"86\u00d759=\", \"65\u00d730=\"],\n  [\"22\u00d776=\", \"88\u00d760=\"],\n  [\"97\u00d744=\", \"65\u00d757=\"],\n  [\"63\u00d754=\", \"82\u00d781=\"],\n  [\"64\u00d733=\", \"87\u00d747=\"],\n  [\"35\u00d744=\", \"24\u00d761=\"],\n  [\"37\u00d725=\", \"82\u00d756=\"],\n  [\"17\u00d771=\", \"85\u00d715=\"],\n  [\"76\u00d737=\", \"33\u00d725=\"],\n  [\"11\u00d756=\", \"12\u00d736=\"],\n  [\"63\u00d799=\", \"89\u00d790=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-by-two-digit multiplication expression in the\n# document with its new value, matched by exact original text (each\n# original expression is unique in the document, so a plain Find/Replace\n# per pair is unambiguous).\n$pairs = @(\n  @(\"26\u00d725=\", \"72\u00d721=\"),\n  @(\"92\u00d798=\", \"62\u00d722=\"),\n  @(\"79\u00d736=\", \"74\u00d748=\"),\n  @(\"16\u00d776=\", \"33\u00d747=\"),\n  @(\"50\u00d793=\", \"65\u00d753=\"),\n  @(\"39\u00d766=\", \"79\u00d754=\"),\n  @(\"69\u00d788=\", \"21\u00d790=\"),\n  @(\"14\u00d725=\", \"22\u00d787=\"),\n  @(\"56\u00d779=\", \"64\u00d774=\"),\n  @(\"53\u00d713=\", \"58\u00d798=\"),\n  @(\"14\u00d786=\", \"61\u00d739=\"),\n  @(\"92\u00d775=\", \"94\u00d758=\"),\n  @(\"66\u00d723=\", \"87\u00d728=\"),\n  @(\"80\u00d742=\", \"36\u00d746=\"),\n  @(\"86\u00d759=\", \"65\u00d730=\"),\n  @(\"22\u00d776=\", \"88\u00d760=\"),\n  @(\"97\u00d744=\", \"65\u00d757=\"),\n  @(\"63\u00d754=\", \"82\u00d781=\"),\n  @(\"64\u00d733=\", \"87\u00d747=\"),\n  @(\"35\u00d744=\", \"24\u00d761=\"),\n  @(\"37\u00d725=\", \"82\u00d756=\"),\n  @(\"17\u00d771=\", \"85\u00d715=\"),\n  @(\"76\u00d737=\", \"33\u00d725=\"),\n  @(\"11\u00d756=\", \"12\u00d736=\"),\n  @(\"63\u00d799=\", \"89\u00d790=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
